$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the poly(A) isolation protocol kit name (shared string text change):
#    "NEBNextPoly(A)E7490" -> "NEBNextPoly(A)E7490L" for every cell that uses it
#    (column G, rows 2:37 - the "polyAIsolationProtocol" column)
$ws.Range("G2:G37").Value = "NEBNextPoly(A)E7490L"

# 2. Widen column G (polyAIsolationProtocol) to fit the longer text, leaving the
#    other columns at their original (default) width.
$ws.Columns.Item(7).ColumnWidth = 24.43

# 3. Update the active selection to the widened column G (G2:G37, active cell G2)
#    instead of the previous I2:I37 selection.
$null = $ws.Range("G2:G37").Select()
